# numero_cliente.xlsx — "detecta campos opcionales y el cliente por numero"
#
# Summary of changes applied:
#  - bd_clientes (sheet2) gains a new "links" column (C) with lookup URLs
#    depending on the "Bases de datos" type (PST vs EFICACIA sources).
#  - A couple of stray formatted-but-empty cells (F12, G15) show up on
#    bd_clientes (underline style carried over from copy/paste of links).
#  - The active sheet switches from bd_clientes back to Hoja1, with the
#    selection left on A6 and column A widened to fit the longer numbers.
#  - bd_clientes' own selection moves to C1 (the new column) and it is no
#    longer the tabSelected sheet.

$wb = $excel.ActiveWorkbook

$hoja1 = $wb.Worksheets.Item("Hoja1")
$bdClientes = $wb.Worksheets.Item("bd_clientes")

# --- bd_clientes: new "links" column -----------------------------------
$dataBase = '"http://192.168.169.23:8083/dbaexperts/dataBase" '
$pstDataBase = '"http://192.168.169.23:8083/dbaexperts/pst_dataBase" '

# Write in the same order the shared-string table records them in the
# target workbook: dataBase first, then pst_dataBase, then the "links"
# header last.
$bdClientes.Range("C5").Value = $dataBase
$bdClientes.Range("C6").Value = $dataBase
$bdClientes.Range("C7").Value = $dataBase

$bdClientes.Range("C2").Value = $pstDataBase
$bdClientes.Range("C3").Value = $pstDataBase
$bdClientes.Range("C4").Value = $pstDataBase

$bdClientes.Range("C1").Value = "links"

# Stray underlined-but-empty cells left over on the sheet.
$bdClientes.Range("F12").Font.Underline = 1
$bdClientes.Range("G15").Font.Underline = 1

# bd_clientes selection moves to C1, and it stops being the active tab.
$null = $bdClientes.Range("C1").Select()

# --- Hoja1: widen column A, move selection, make it the active tab -----
$hoja1.Columns.Item(1).ColumnWidth = 31.608072916666668
$null = $hoja1.Range("A6").Select()
$null = $hoja1.Activate()
